$d = $word.ActiveDocument

function Set-ParaXml($paraIndex, $innerXml) {
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# --- Paragraph 53: "Anh Lê Văn Long: " -> "Anh Trần Văn Cường " (drop italic, new name) ---
Set-ParaXml 53 '<w:body><w:p><w:pPr><w:rPr><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve">Anh </w:t></w:r><w:r><w:rPr><w:iCs/></w:rPr><w:t>Trần Văn Cường</w:t></w:r><w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body>'

# --- Paragraph 55: "Lập trình viên:  Pham Lan" -> "Lập " / "trình viên: " / "Lê Quang Dũng" (drop italic) ---
Set-ParaXml 55 '<w:body><w:p><w:pPr><w:rPr><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve">Lập </w:t></w:r><w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve">trình viên: </w:t></w:r><w:r><w:rPr><w:iCs/></w:rPr><w:t>Lê Quang Dũng</w:t></w:r></w:p></w:body>'

# --- Paragraph 56: "Phiên dịch: Ngọc, " -> "Phiên dịch: " / "Trần Quang Anh" (drop italic) ---
Set-ParaXml 56 '<w:body><w:p><w:pPr><w:rPr><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve">Phiên dịch: </w:t></w:r><w:r><w:rPr><w:iCs/></w:rPr><w:t>Trần Quang Anh</w:t></w:r></w:p></w:body>'

# Paragraph 57 (Heading2 "Phân chia vai trò...") is unchanged in content; its bookmark id
# renumbers automatically once the _GoBack bookmark is relocated below.

# --- Paragraph 58: "Giám đốc: Hưng, ..." -> "Giám đốc: " / "Hoàng Công Hậu" (drop italic) ---
Set-ParaXml 58 '<w:body><w:p><w:pPr><w:rPr><w:iCs/></w:rPr></w:pPr><w:r><w:rPr><w:iCs/></w:rPr><w:t xml:space="preserve">Giám đốc: </w:t></w:r><w:r><w:rPr><w:iCs/></w:rPr><w:t>Hoàng Công Hậu</w:t></w:r></w:p></w:body>'

# --- Paragraphs 59-60: "Trung: IT..." + "Phiên dịch: Bích" -> single paragraph
#     "Trao đổi với khách hàng: Đoàn Đức Bảo" (no pPr, drop italic) ---
$rng59 = $d.Paragraphs(59).Range
$rng60 = $d.Paragraphs(60).Range
$combined = $d.Range($rng59.Start, $rng60.End)
$pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:iCs/></w:rPr><w:t>Trao đổi với khách hàng: Đoàn Đức Bảo</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$combined.InsertXML($pkg)

# --- Relocate the _GoBack bookmark: remove the old one, add a new one between
#     "trình viên: " and "Lê Quang Dũng" in paragraph 55 ---
$old = $d.Bookmarks("_GoBack")
$old.Delete()

$p55 = $d.Paragraphs(55)
$insPoint = $d.Range($p55.Range.Start, $p55.Range.Start)
$insPoint.Find.Execute("trình viên: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$full55 = $p55.Range
$full55.Find.Execute("trình viên: ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmPoint = $d.Range($full55.End, $full55.End)
$d.Bookmarks.Add("_GoBack", $bmPoint)

# --- Remove one of the three empty paragraphs following (3 -> 2) ---
$i = 0
$emptyIdx = @()
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($i -ge 59 -and $i -le 65) {
        if ($p.Range.Text -eq "`r") {
            $emptyIdx += $i
        }
    }
}
$lastEmpty = $emptyIdx[$emptyIdx.Length - 1]
$d.Paragraphs($lastEmpty).Range.Delete()

Write-Output "done"
